# Grade entries for rows 22 and 24 (S.No 14 "whoPurchasedProduct() method" and
# S.No 16 "findAllBrands()") in the CustomerMapping Class table, plus moving the
# active selection/view to the area just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: points earned (E22) and a grading comment (F22) explaining the deduction.
$ws.Range("E22").Value = 6
$ws.Range("F22").Value = " -4 for wrong logic in if condition"

# Row 24: full points earned.
$ws.Range("E24").Value = 10

# Move the selection to reflect where grading left off.
$ws.Range("E24").Select()
